$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update most_frequent_value (Q4)
$ws.Range("Q4").Value = "Cannon Hill Station"

# Row 5: update most_frequent_value (Q5)
$ws.Range("Q5").Value = "Grey Street, South Brisbane, Brisbane Central"

# Row 6: update blank_count (P6), most_frequent_value (Q6),
# memory_consumed_bytes (R6), pattern_count (S6), patterns (T6)
$ws.Range("P6").Value = 3
$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "-27.4818"
$ws.Range("R6").Value = 7787
$ws.Range("S6").Value = 756
$ws.Range("T6").Value = "['-26.1627', '-26.366181', '-26.386766', '-26.386933', '-26.387332', '-26.3921', '-26.39659', '-26.3971', '-26.398385', '-26.3986', '-26.399192', '-26.4014', '-26.4158', '-26.4171', '-26.4177', '-26.4184', '-26.4187', '-26.475572', '-26.4812', '-26.529825', '-26.5309', '-26.559987', '-26.6007', '-26.61413', '-26.616492', '-26.6221', '-26.6258', '-26.6278', '-26.6547', '-26.654875', '-26.654877', '-26.6549', '-26.6605', '-26.6608', '-26.6624', '-26.66553', '-26.668315', '-26.680953', '-26.682702', '-26.6841', '-26.687894', '-26.6886', '-26.70291', '-26.703187', '-26.7152', '-26.718071', '-26.718545', '-26.7197', '-26.7259', '-26.759023', '-26.7653', '-26.773831', '-26.7899', '-26.80477', '-26.804983', '-26.8084', '-26.808765', '-26.821495', '-26.857889', '-26.857894', '-26.8997', '-26.959839', '-27.0113', '-27.061943', '-27.064302', '-27.06552', '-27.0683', '-27.076', '-27.0833', '-27.083338', '-27.083347', '-27.1098', '-27.1123', '-27.154338', '-27.163063', '-27.163078', '-27.173203', '-27.19383', '-27.196441', '-27.202777', '-27.2029', '-27.2173', '-27.220542', '-27.225499', '-27.2256', '-27.227944', '-27.2285', '-27.2305', '-27.230745', '-27.232057', '-27.238184', '-27.238653', '-27.2398', '-27.244551', '-27.2501', '-27.250229', '-27.2609', '-27.2678', '-27.2732', '-27.2816', '-27.2915', '-27.2971', '-27.299357', '-27.299367', '-27.3069', '-27.307691', '-27.310447', '-27.3118', '-27.32111', '-27.321362', '-27.321372', '-27.3222', '-27.324331', '-27.3269', '-27.3273', '-27.3284', '-27.330219', '-27.330283', '-27.3307', '-27.3324', '-27.337817', '-27.339', '-27.342044', '-27.3453', '-27.346377', '-27.346405', '-27.3482', '-27.348448', '-27.349545', '-27.3507', '-27.3527', '-27.3581', '-27.358375', '-27.3593', '-27.3595', '-27.3626', '-27.363933', '-27.363956', '-27.364', '-27.3656', '-27.3689', '-27.372756', '-27.3749', '-27.374986', '-27.374995', '-27.3753', '-27.379966', '-27.3829', '-27.382901', '-27.384486', '-27.3852', '-27.3854', '-27.3893', '-27.3903', '-27.3917', '-27.3936', '-27.3973', '-27.4003', '-27.4008', '-27.4009', '-27.4013', '-27.4028', '-27.403', '-27.4042', '-27.4045', '-27.405', '-27.406636', '-27.4073', '-27.408392', '-27.4086', '-27.409', '-27.409867', '-27.410027', '-27.4115', '-27.4121', '-27.415092', '-27.4151', '-27.4175', '-27.4193', '-27.42', '-27.4203', '-27.422', '-27.4222', '-27.422296', '-27.422652', '-27.423433', '-27.423443', '-27.4247', '-27.424959', '-27.4266', '-27.429214', '-27.429409', '-27.429438', '-27.4314', '-27.4322', '-27.4328', '-27.4352', '-27.436058', '-27.4362', '-27.4368', '-27.4371', '-27.4372', '-27.438428', '-27.438633', '-27.4397', '-27.441', '-27.4421', '-27.44285', '-27.4447', '-27.4449', '-27.445', '-27.4451', '-27.445162', '-27.445181', '-27.4452', '-27.4455', '-27.445599', '-27.445752', '-27.446136', '-27.4462', '-27.4474', '-27.4475', '-27.449758', '-27.450247', '-27.451', '-27.4518', '-27.452', '-27.4533', '-27.4536', '-27.454', '-27.454744', '-27.4554', '-27.455465', '-27.4567', '-27.456767', '-27.45679', '-27.4568', '-27.456839', '-27.457', '-27.4571', '-27.4574', '-27.457771', '-27.458', '-27.4589', '-27.458927', '-27.459', '-27.4601', '-27.4611', '-27.4616', '-27.462', '-27.4632', '-27.4646', '-27.4648', '-27.4649', '-27.4652', '-27.4653', '-27.465367', '-27.465376', '-27.465434', '-27.465905', '-27.466', '-27.4661', '-27.466104', '-27.4662', '-27.4664', '-27.4668', '-27.4669', '-27.466974', '-27.4671', '-27.467113', '-27.467361', '-27.4674', '-27.4675', '-27.468', '-27.4683', '-27.4685', '-27.4686', '-27.4687', '-27.468872', '-27.4689', '-27.4693', '-27.4694', '-27.469414', '-27.469461', '-27.46957', '-27.4697', '-27.4701', '-27.4702', '-27.4703', '-27.4704', '-27.4705', '-27.470582', '-27.4706', '-27.470755', '-27.4709', '-27.471', '-27.4711', '-27.4712', '-27.4713', '-27.4714', '-27.4715', '-27.4721', '-27.472178', '-27.472558', '-27.4732', '-27.473771', '-27.4744', '-27.474857', '-27.475372', '-27.4755', '-27.4757', '-27.4759', '-27.476547', '-27.476557', '-27.4769', '-27.477569', '-27.4778', '-27.4781', '-27.4785', '-27.4786', '-27.47948', '-27.47966', '-27.4797', '-27.479865', '-27.4807', '-27.48086', '-27.4814', '-27.4818', '-27.4827', '-27.483568', '-27.4842', '-27.484423', '-27.4845', '-27.4847', '-27.4848', '-27.48533', '-27.485334', '-27.4854', '-27.4857', '-27.486', '-27.4864', '-27.486614', '-27.487', '-27.487198', '-27.4875', '-27.4877', '-27.488', '-27.4897', '-27.4905', '-27.491928', '-27.492742', '-27.492889', '-27.493', '-27.4936', '-27.4937', '-27.493933', '-27.493943', '-27.494284', '-27.4944', '-27.4955', '-27.4963', '-27.496571', '-27.496576', '-27.496863', '-27.497029', '-27.497168', '-27.497447', '-27.4976', '-27.497963', '-27.498', '-27.498036', '-27.4982', '-27.498635', '-27.499', '-27.4992', '-27.499652', '-27.4998', '-27.500757', '-27.5025', '-27.5026', '-27.502998', '-27.503', '-27.503265', '-27.503369', '-27.504137', '-27.504502', '-27.505594', '-27.506749', '-27.5072', '-27.5074', '-27.5085', '-27.509', '-27.5095', '-27.5102', '-27.510946', '-27.5112', '-27.511648', '-27.5121', '-27.5146', '-27.5163', '-27.5178', '-27.5184', '-27.5195', '-27.519645', '-27.5199', '-27.5213', '-27.521484', '-27.5217', '-27.5219', '-27.5231', '-27.523145', '-27.5239', '-27.523933', '-27.524102', '-27.52421', '-27.524219', '-27.5246', '-27.5253', '-27.5255', '-27.5258', '-27.5275', '-27.5311', '-27.5316', '-27.531895', '-27.5324', '-27.5326', '-27.5341', '-27.534681', '-27.535', '-27.53513', '-27.5355', '-27.5373', '-27.5388', '-27.5396', '-27.543033', '-27.543289', '-27.5447', '-27.545804', '-27.5459', '-27.5466', '-27.547', '-27.5481', '-27.5483', '-27.54897', '-27.549068', '-27.550683', '-27.5525', '-27.5527', '-27.5529', '-27.5532', '-27.555213', '-27.55529', '-27.555737', '-27.5559', '-27.555952', '-27.558271', '-27.5612', '-27.561953', '-27.5622', '-27.5631', '-27.5634', '-27.5639', '-27.5641', '-27.5659', '-27.566231', '-27.567043', '-27.5684', '-27.569609', '-27.570008', '-27.5706', '-27.571407', '-27.572561', '-27.574419', '-27.57471', '-27.575481', '-27.576019', '-27.576024', '-27.57848', '-27.5792', '-27.5808', '-27.582671', '-27.58269', '-27.584466', '-27.5854', '-27.5865', '-27.5874', '-27.589', '-27.590132', '-27.5902', '-27.5915', '-27.5924', '-27.5938', '-27.595', '-27.595882', '-27.5964', '-27.596442', '-27.5979', '-27.597943', '-27.5983', '-27.59889', '-27.599', '-27.5995', '-27.6014', '-27.6018', '-27.6033', '-27.6045', '-27.6047', '-27.6049', '-27.6059', '-27.6061', '-27.607', '-27.6073', '-27.6075', '-27.607519', '-27.608385', '-27.6088', '-27.609357', '-27.60946', '-27.61', '-27.610074', '-27.6109', '-27.6112', '-27.6125', '-27.6129', '-27.6135', '-27.6144', '-27.614579', '-27.6147', '-27.6148', '-27.61511', '-27.615168', '-27.6187', '-27.620445', '-27.624', '-27.6255', '-27.6266', '-27.628126', '-27.6284', '-27.629638', '-27.630829', '-27.6314', '-27.6334', '-27.634837', '-27.637003', '-27.6382', '-27.6413', '-27.6419', '-27.6427', '-27.6429', '-27.642944', '-27.6432', '-27.6449', '-27.645243', '-27.648638', '-27.648659', '-27.651875', '-27.6533', '-27.65422', '-27.6547', '-27.6574', '-27.660196', '-27.660298', '-27.6603', '-27.660306', '-27.660548', '-27.660573', '-27.6606', '-27.662644', '-27.6639', '-27.664117', '-27.664398', '-27.664773', '-27.666518', '-27.671', '-27.6734', '-27.6735', '-27.674276', '-27.674482', '-27.6776', '-27.678743', '-27.680633', '-27.688452', '-27.688471', '-27.690521', '-27.6924', '-27.696258', '-27.6974', '-27.6982', '-27.7008', '-27.7078', '-27.7095', '-27.7134', '-27.7162', '-27.717115', '-27.717125', '-27.7383', '-27.7444', '-27.7639', '-27.78611', '-27.7991', '-27.830644', '-27.84369', '-27.8535', '-27.868739', '-27.878122', '-27.878307', '-27.880533', '-27.885812', '-27.912639', '-27.913145', '-27.923565', '-27.9258', '-27.926', '-27.926171', '-27.9272', '-27.931162', '-27.931558', '-27.931705', '-27.946272', '-27.947034', '-27.9487', '-27.959043', '-27.96055', '-27.96094', '-27.962054', '-27.962918', '-27.962972', '-27.963057', '-27.966594', '-27.967446', '-27.967624', '-27.968746', '-27.969082', '-27.9695', '-27.970143', '-27.970411', '-27.970573', '-27.970686', '-27.971077', '-27.971204', '-27.973236', '-27.97349', '-27.973494', '-27.973923', '-27.977936', '-27.979071', '-27.979344', '-27.979743', '-27.981359', '-27.9821', '-27.982328', '-27.984646', '-27.985951', '-27.9871', '-27.98768', '-27.9891', '-27.989739', '-27.99073', '-27.991322', '-27.993062', '-27.993107', '-27.9933', '-27.993316', '-27.9936', '-27.9943', '-27.995222', '-27.995515', '-27.995758', '-27.996317', '-27.996461', '-27.997265', '-27.997786', '-27.997806', '-27.998069', '-27.998162', '-27.999044', '-27.999085', '-27.999151', '-27.999278', '-27.999348', '-27.999546', '-28.000986', '-28.001165', '-28.001446', '-28.001657', '-28.001777', '-28.001927', '-28.00215', '-28.002304', '-28.003255', '-28.003694', '-28.005486', '-28.005792', '-28.006013', '-28.006222', '-28.006307', '-28.00749', '-28.009103', '-28.009639', '-28.010543', '-28.011257', '-28.012158', '-28.013906', '-28.015495', '-28.016302', '-28.016508', '-28.017705', '-28.024987', '-28.026559', '-28.027816', '-28.02804', '-28.028599', '-28.029247', '-28.029567', '-28.030702', '-28.031868', '-28.033787', '-28.035402', '-28.035645', '-28.035763', '-28.035828', '-28.036277', '-28.040189', '-28.045714', '-28.0465', '-28.049584', '-28.0587', '-28.065476', '-28.0714', '-28.077591', '-28.078462', '-28.080915', '-28.084951', '-28.085385', '-28.088982', '-28.090215', '-28.094794', '-28.096499', '-28.096893', '-28.109083', '-28.118114', '-28.119505', '-28.1223', '-28.135343', '-28.1355', '-28.14411', '-28.165645', '-28.16626', '-28.1678', '-28.168261', '-28.16833', '-28.169044', '-28.171', '-28.171647', '-28.207807', '0']"

# Row 7: update blank_count (P7), memory_consumed_bytes (R7),
# pattern_count (S7), patterns (T7)
$ws.Range("P7").Value = 3
$ws.Range("R7").Value = 7329
$ws.Range("S7").Value = 603
$ws.Range("T7").Value = "['0', '151.930743', '151.954786', '151.968969', '152.270886', '152.277841', '152.279529', '152.337285', '152.593', '152.63', '152.665', '152.668', '152.686', '152.709', '152.732', '152.746', '152.751467', '152.754', '152.756307', '152.759836', '152.761', '152.774425', '152.775', '152.789748', '152.789791', '152.789796', '152.79', '152.804549', '152.804571', '152.822', '152.832', '152.833', '152.851', '152.85505', '152.868', '152.871', '152.872', '152.872206', '152.886495', '152.886854', '152.897662', '152.902', '152.902834', '152.91', '152.911', '152.912', '152.918', '152.919', '152.921474', '152.922', '152.923094', '152.923292', '152.924', '152.929', '152.930148', '152.935', '152.935247', '152.935947', '152.936888', '152.938', '152.938354', '152.939', '152.94', '152.943656', '152.946', '152.946558', '152.948', '152.948024', '152.949', '152.949775', '152.95', '152.950648', '152.951', '152.952', '152.953', '152.954011', '152.954021', '152.956586', '152.957', '152.957107', '152.957113', '152.958', '152.958349', '152.958386', '152.959', '152.959021', '152.959402', '152.96021', '152.960821', '152.961', '152.962', '152.962742', '152.963', '152.964943', '152.966', '152.968', '152.969', '152.969461', '152.972679', '152.972871', '152.972882', '152.973', '152.974', '152.97408', '152.97469', '152.975', '152.976', '152.978784', '152.9788', '152.979294', '152.98', '152.981', '152.982', '152.983', '152.983028', '152.986', '152.986032', '152.986145', '152.986997', '152.987', '152.99', '152.990315', '152.990456', '152.990995', '152.991', '152.99154', '152.992', '152.992052', '152.993', '152.993461', '152.994', '152.995942', '152.996732', '152.996743', '152.998', '152.999', '152.999797', '153.000265', '153.000987', '153.003', '153.005', '153.006', '153.006742', '153.007', '153.007559', '153.010862', '153.011', '153.011644', '153.011654', '153.01192', '153.012', '153.012848', '153.013', '153.014', '153.014624', '153.015', '153.01509', '153.015968', '153.017', '153.017093', '153.017253', '153.017324', '153.017466', '153.018', '153.018268', '153.019', '153.019897', '153.02', '153.020663', '153.021', '153.021132', '153.021835', '153.022', '153.022332', '153.023', '153.02326', '153.02343', '153.023727', '153.024', '153.025', '153.025012', '153.025228', '153.025309', '153.025936', '153.025944', '153.026', '153.026006', '153.026505', '153.026526', '153.026739', '153.027', '153.027863', '153.028', '153.029', '153.029145', '153.029161', '153.029327', '153.03', '153.030513', '153.030874', '153.031', '153.032', '153.032021', '153.032094', '153.033', '153.03302', '153.033095', '153.033842', '153.033977', '153.034', '153.034598', '153.035', '153.035526', '153.036', '153.036245', '153.036824', '153.037', '153.037549', '153.03756', '153.038', '153.038964', '153.039', '153.039182', '153.039488', '153.04', '153.040021', '153.040386', '153.040397', '153.040452', '153.040492', '153.040493', '153.040967', '153.041', '153.041259', '153.041343', '153.04178', '153.042', '153.043', '153.044', '153.044459', '153.04446', '153.045', '153.0452', '153.045734', '153.046', '153.046161', '153.046177', '153.046272', '153.047', '153.04764', '153.048923', '153.049', '153.05', '153.050149', '153.050647', '153.051', '153.052821', '153.053', '153.054', '153.054055', '153.054448', '153.054915', '153.055', '153.056', '153.057', '153.057646', '153.058', '153.059', '153.06', '153.060114', '153.061', '153.061427', '153.061561', '153.061564', '153.061919', '153.062', '153.062414', '153.063', '153.063238', '153.063376', '153.063397', '153.063761', '153.06389', '153.064', '153.064185', '153.064402', '153.064916', '153.065', '153.065128', '153.065223', '153.066', '153.066338', '153.067622', '153.067633', '153.067783', '153.068', '153.068175', '153.069', '153.07', '153.071304', '153.07245', '153.074', '153.074011', '153.075', '153.075345', '153.077', '153.077268', '153.077279', '153.079', '153.079498', '153.08', '153.080204', '153.081', '153.082', '153.08257', '153.083', '153.083098', '153.083234', '153.084', '153.084894', '153.087227', '153.088', '153.088099', '153.088115', '153.089', '153.089098', '153.089572', '153.089624', '153.09', '153.090593', '153.090858', '153.091', '153.091151', '153.091411', '153.091521', '153.091644', '153.092', '153.093399', '153.09351', '153.095', '153.096', '153.096887', '153.098', '153.1', '153.101', '153.101627', '153.101739', '153.10215', '153.102971', '153.104', '153.105', '153.105607', '153.105673', '153.106435', '153.106621', '153.107', '153.107871', '153.108', '153.109', '153.11', '153.112', '153.114', '153.11427', '153.11473', '153.115019', '153.115212', '153.115424', '153.117566', '153.11775', '153.118', '153.119', '153.120058', '153.120477', '153.123', '153.124', '153.124299', '153.124547', '153.124881', '153.125', '153.125012', '153.125323', '153.126', '153.126052', '153.127', '153.128403', '153.129027', '153.13', '153.131188', '153.135425', '153.136', '153.139', '153.142', '153.145', '153.148047', '153.153', '153.153712', '153.154778', '153.156', '153.157087', '153.158635', '153.158762', '153.159', '153.166', '153.169', '153.17', '153.170772', '153.171', '153.172', '153.173', '153.173541', '153.174', '153.177819', '153.18', '153.185', '153.185659', '153.187', '153.188', '153.191099', '153.193', '153.194', '153.201', '153.203', '153.204', '153.204215', '153.204236', '153.206', '153.211', '153.216', '153.218472', '153.218483', '153.221', '153.237809', '153.238894', '153.241914', '153.245', '153.256', '153.258136', '153.265657', '153.265667', '153.267', '153.280957', '153.281', '153.282', '153.282808', '153.283', '153.298', '153.301', '153.308834', '153.313517', '153.314545', '153.318', '153.321239', '153.323605', '153.323755', '153.32816', '153.332', '153.333025', '153.334068', '153.336', '153.339', '153.35', '153.350839', '153.364813', '153.365036', '153.365587', '153.366608', '153.377048', '153.378', '153.379559', '153.379764', '153.38047', '153.380514', '153.380708', '153.380854', '153.381172', '153.381709', '153.382003', '153.382255', '153.383257', '153.384389', '153.385406', '153.387762', '153.387871', '153.390741', '153.391622', '153.393909', '153.399049', '153.399531', '153.39985', '153.401878', '153.403764', '153.403905', '153.407014', '153.407677', '153.408', '153.408303', '153.409', '153.410518', '153.412273', '153.412793', '153.413301', '153.414', '153.414091', '153.414237', '153.414436', '153.414542', '153.414785', '153.415458', '153.415665', '153.415717', '153.415815', '153.418525', '153.418964', '153.419628', '153.420625', '153.422847', '153.423689', '153.423854', '153.423886', '153.42483', '153.424967', '153.425076', '153.426261', '153.426284', '153.426351', '153.426385', '153.426632', '153.426897', '153.426998', '153.42705', '153.427522', '153.427723', '153.427801', '153.427978', '153.428026', '153.428313', '153.428591', '153.428659', '153.428716', '153.428721', '153.428729', '153.42889', '153.42891', '153.428923', '153.428966', '153.428985', '153.429105', '153.429112', '153.429242', '153.429275', '153.429349', '153.429371', '153.429379', '153.429393', '153.42946', '153.429523', '153.429638', '153.429727', '153.429731', '153.429787', '153.42991', '153.429921', '153.430115', '153.430342', '153.430424', '153.430453', '153.430599', '153.430639', '153.430727', '153.43078', '153.43085', '153.430861', '153.431112', '153.431211', '153.431294', '153.431904', '153.431953', '153.432548', '153.432791', '153.434403', '153.435431', '153.436301', '153.438', '153.442285', '153.446451', '153.448477', '153.451783', '153.465528', '153.469333', '153.47051', '153.471363', '153.472', '153.472662', '153.489', '153.509074', '153.512886', '153.528', '153.536129', '153.541951', '153.542848', '153.542908', '153.543691', '153.545211']"
